# Add new power plants to Electricity Source subscript (issues #280 and #99)
# Appends six new plant types to the "PTUfIGaMDC" worksheet, each with a
# default flag value of 0 in column B, just like the existing rows.

$wb = $excel.ActiveWorkbook

$aboutSheet = $wb.Worksheets.Item("About")
$ws = $wb.Worksheets.Item("PTUfIGaMDC")

$newPlantTypes = @(
    "hard coal w CCS",
    "natural gas combined cycle w CCS",
    "biomass w CCS",
    "lignite w CCS",
    "small modular reactor",
    "hydrogen"
)

$row = 19
foreach ($plantType in $newPlantTypes) {
    $ws.Cells.Item($row, 1).Value = $plantType
    $ws.Cells.Item($row, 2).Value = 0
    $row++
}

# Update the selection on the PTUfIGaMDC sheet to the cell just below the
# newly-added data (matches the cursor position left behind after typing
# in the new rows), then restore "About" as the active sheet/tab.
$ws.Activate()
$ws.Range("A25").Select()
$aboutSheet.Activate()
